$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row (71) of prediction data to the bottom of the table.
$row = 71
$ws.Cells.Item($row, 1).Value = 0.333333
$ws.Cells.Item($row, 2).Value = 0.111111
$ws.Cells.Item($row, 3).Value = 0.333333
$ws.Cells.Item($row, 4).Value = 0.888888
$ws.Cells.Item($row, 5).Value = 0.222222
$ws.Cells.Item($row, 6).Value = 0.555555
$ws.Cells.Item($row, 7).Value = 0.07444996695619652
$ws.Cells.Item($row, 8).Value = "query"
